$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = [double]"21.13602999991303"
$ws.Cells.Item(2, 3).Value = [double]"381"
$ws.Cells.Item(2, 5).Value = [double]"-7.720553435036679E-07"
$ws.Cells.Item(2, 6).Value = [double]"0.2429359398155243"
$ws.Cells.Item(2, 7).Value = [double]"3657.212172385902"
$ws.Cells.Item(2, 8).Value = [double]"0.5779273666292173"
$ws.Cells.Item(3, 2).Value = [double]"21.87499382997181"
$ws.Cells.Item(3, 3).Value = [double]"5"
$ws.Cells.Item(3, 6).Value = [double]"0.2810465917672703"
$ws.Cells.Item(3, 7).Value = [double]"3754.341545634866"
$ws.Cells.Item(3, 8).Value = [double]"0.5826585984273498"
$ws.Cells.Item(4, 2).Value = [double]"22.61379504996888"
$ws.Cells.Item(4, 5).Value = [double]"3.860263435036679E-07"
$ws.Cells.Item(4, 6).Value = [double]"0.3023404908804053"
$ws.Cells.Item(4, 7).Value = [double]"3926.549510483434"
$ws.Cells.Item(4, 8).Value = [double]"0.5759202829250633"
$ws.Cells.Item(5, 2).Value = [double]"23.31961550990591"
$ws.Cells.Item(5, 6).Value = [double]"0.3155688877540641"
$ws.Cells.Item(5, 7).Value = [double]"4108.00129283187"
$ws.Cells.Item(5, 8).Value = [double]"0.567663295301215"
$ws.Cells.Item(6, 6).Value = [double]"0.3223510320401817"
$ws.Cells.Item(6, 7).Value = [double]"4297.407113996655"
$ws.Cells.Item(6, 8).Value = [double]"0.5609274348581729"
$ws.Cells.Item(7, 2).Value = [double]"24.99199342987444"
$ws.Cells.Item(7, 4).Value = [double]"7"
$ws.Cells.Item(7, 5).Value = [double]"50.12566026499999"
$ws.Cells.Item(7, 6).Value = [double]"0.3663702942492986"
$ws.Cells.Item(7, 7).Value = [double]"4507.338960230962"
$ws.Cells.Item(7, 8).Value = [double]"0.5544733522458187"
$ws.Cells.Item(8, 2).Value = [double]"25.98912602995889"
$ws.Cells.Item(8, 5).Value = [double]"50.125660265"
$ws.Cells.Item(8, 6).Value = [double]"0.3660677534359533"
$ws.Cells.Item(8, 7).Value = [double]"4664.08325602357"
$ws.Cells.Item(8, 8).Value = [double]"0.5572183128676885"
$ws.Cells.Item(9, 2).Value = [double]"27.22151567996369"
$ws.Cells.Item(9, 4).Value = [double]"10"
$ws.Cells.Item(9, 5).Value = [double]"67.38647562243641"
$ws.Cells.Item(9, 6).Value = [double]"0.3854132385370957"
$ws.Cells.Item(9, 7).Value = [double]"4805.90949083478"
$ws.Cells.Item(9, 8).Value = [double]"0.5664175684514473"
$ws.Cells.Item(10, 2).Value = [double]"27.85923207996592"
$ws.Cells.Item(10, 4).Value = [double]"18"
$ws.Cells.Item(10, 5).Value = [double]"220.2353592769195"
$ws.Cells.Item(10, 6).Value = [double]"0.3828738558137177"
$ws.Cells.Item(10, 7).Value = [double]"4987.251217864619"
$ws.Cells.Item(10, 8).Value = [double]"0.5586089583811732"
$ws.Cells.Item(11, 2).Value = [double]"28.46385891996019"
$ws.Cells.Item(11, 3).Value = [double]"0"
$ws.Cells.Item(11, 4).Value = [double]"18"
$ws.Cells.Item(11, 5).Value = [double]"190.6254554704508"
$ws.Cells.Item(11, 6).Value = [double]"0.3774139506907817"
$ws.Cells.Item(11, 7).Value = [double]"5195.005815058338"
$ws.Cells.Item(11, 8).Value = [double]"0.5479081243269128"
$ws.Cells.Item(12, 2).Value = [double]"29.10164526999474"
$ws.Cells.Item(12, 3).Value = [double]"15"
$ws.Cells.Item(12, 4).Value = [double]"30"
$ws.Cells.Item(12, 5).Value = [double]"364.3634029514913"
$ws.Cells.Item(12, 6).Value = [double]"0.3310637573392921"
$ws.Cells.Item(12, 7).Value = [double]"5468.825389380543"
$ws.Cells.Item(12, 8).Value = [double]"0.5321370348833007"
$ws.Cells.Item(13, 2).Value = [double]"29.74282199996686"
$ws.Cells.Item(13, 3).Value = [double]"5.000000013992965"
$ws.Cells.Item(13, 5).Value = [double]"360.1078749434695"
$ws.Cells.Item(13, 6).Value = [double]"0.3199459227193345"
$ws.Cells.Item(13, 7).Value = [double]"5712.71101711458"
$ws.Cells.Item(13, 8).Value = [double]"0.5206428595961011"
$ws.Cells.Item(14, 2).Value = [double]"30.3389967399616"
$ws.Cells.Item(14, 4).Value = [double]"37"
$ws.Cells.Item(14, 5).Value = [double]"263.3234432383766"
$ws.Cells.Item(14, 6).Value = [double]"0.3225882980186723"
$ws.Cells.Item(14, 7).Value = [double]"5911.05697197583"
$ws.Cells.Item(14, 8).Value = [double]"0.5132584051176974"
$ws.Cells.Item(15, 2).Value = [double]"30.55337515996858"
$ws.Cells.Item(15, 3).Value = [double]"6.999999986158552"
$ws.Cells.Item(15, 4).Value = [double]"38"
$ws.Cells.Item(15, 5).Value = [double]"329.7690017152662"
$ws.Cells.Item(15, 6).Value = [double]"0.3113104022545897"
$ws.Cells.Item(15, 7).Value = [double]"6020.482564941773"
$ws.Cells.Item(15, 8).Value = [double]"0.5074904682539194"
$ws.Cells.Item(16, 2).Value = [double]"30.78208348995901"
$ws.Cells.Item(16, 3).Value = [double]"12"
$ws.Cells.Item(16, 4).Value = [double]"31"
$ws.Cells.Item(16, 5).Value = [double]"267.012032266509"
$ws.Cells.Item(16, 6).Value = [double]"0.2860594173261237"
$ws.Cells.Item(16, 7).Value = [double]"6159.884177319526"
$ws.Cells.Item(16, 8).Value = [double]"0.4997185434638129"
$ws.Cells.Item(17, 2).Value = [double]"31.00310876998401"
$ws.Cells.Item(17, 3).Value = [double]"5"
$ws.Cells.Item(17, 4).Value = [double]"29"
$ws.Cells.Item(17, 5).Value = [double]"284.4107743913463"
$ws.Cells.Item(17, 6).Value = [double]"0.275847131739776"
$ws.Cells.Item(17, 7).Value = [double]"6274.332603380203"
$ws.Cells.Item(17, 8).Value = [double]"0.4941260008001734"
$ws.Cells.Item(18, 2).Value = [double]"31.19056711995681"
$ws.Cells.Item(18, 3).Value = [double]"2"
$ws.Cells.Item(18, 4).Value = [double]"23"
$ws.Cells.Item(18, 5).Value = [double]"242.9058946119717"
$ws.Cells.Item(18, 6).Value = [double]"0.276403726832346"
$ws.Cells.Item(18, 7).Value = [double]"6431.78682902958"
$ws.Cells.Item(18, 8).Value = [double]"0.4849440435304788"
$ws.Cells.Item(19, 2).Value = [double]"31.36134633995734"
$ws.Cells.Item(19, 3).Value = [double]"5"
$ws.Cells.Item(19, 4).Value = [double]"24"
$ws.Cells.Item(19, 5).Value = [double]"205.9909204948982"
$ws.Cells.Item(19, 6).Value = [double]"0.2762300442159518"
$ws.Cells.Item(19, 7).Value = [double]"6539.650957241288"
$ws.Cells.Item(19, 8).Value = [double]"0.4795568837696336"
$ws.Cells.Item(20, 2).Value = [double]"31.40993452997044"
$ws.Cells.Item(20, 5).Value = [double]"155.1575451265325"
$ws.Cells.Item(20, 6).Value = [double]"0.2729897152345186"
$ws.Cells.Item(20, 7).Value = [double]"6594.367667134434"
$ws.Cells.Item(20, 8).Value = [double]"0.4763145780680978"
$ws.Cells.Item(21, 2).Value = [double]"31.44015586997127"
$ws.Cells.Item(21, 3).Value = [double]"3"
$ws.Cells.Item(21, 5).Value = [double]"83.67055736978689"
$ws.Cells.Item(21, 6).Value = [double]"0.2516711655094153"
$ws.Cells.Item(21, 7).Value = [double]"6633.012362057476"
$ws.Cells.Item(21, 8).Value = [double]"0.4739951345457607"
$ws.Cells.Item(22, 2).Value = [double]"31.4426515899704"
$ws.Cells.Item(22, 4).Value = [double]"15"
$ws.Cells.Item(22, 5).Value = [double]"104.2423199513168"
$ws.Cells.Item(22, 6).Value = [double]"0.2428704594071747"
$ws.Cells.Item(22, 7).Value = [double]"6681.444901230538"
$ws.Cells.Item(22, 8).Value = [double]"0.4705965858399809"
$ws.Cells.Item(23, 2).Value = [double]"31.42525142997249"
$ws.Cells.Item(23, 3).Value = [double]"8.00000005738662"
$ws.Cells.Item(23, 4).Value = [double]"24"
$ws.Cells.Item(23, 5).Value = [double]"180.0298752901483"
$ws.Cells.Item(23, 6).Value = [double]"0.2477134225436995"
$ws.Cells.Item(23, 7).Value = [double]"6822.208081542055"
$ws.Cells.Item(23, 8).Value = [double]"0.460631676055083"
$ws.Cells.Item(24, 2).Value = [double]"31.37385642997133"
$ws.Cells.Item(24, 3).Value = [double]"3"
$ws.Cells.Item(24, 4).Value = [double]"20"
$ws.Cells.Item(24, 5).Value = [double]"140.8887298231378"
$ws.Cells.Item(24, 6).Value = [double]"0.250105382835819"
$ws.Cells.Item(24, 7).Value = [double]"6880.294698091778"
$ws.Cells.Item(24, 8).Value = [double]"0.4559958229503272"
$ws.Cells.Item(25, 2).Value = [double]"31.22650743008037"
$ws.Cells.Item(25, 3).Value = [double]"4.999999972800943"
$ws.Cells.Item(25, 4).Value = [double]"21"
$ws.Cells.Item(25, 5).Value = [double]"125.8153796792718"
$ws.Cells.Item(25, 6).Value = [double]"0.2467965335534115"
$ws.Cells.Item(25, 7).Value = [double]"6873.929182424129"
$ws.Cells.Item(25, 8).Value = [double]"0.4542745000912007"
$ws.Cells.Item(26, 2).Value = [double]"31.07599652997714"
$ws.Cells.Item(26, 3).Value = [double]"12"
$ws.Cells.Item(26, 4).Value = [double]"30"
$ws.Cells.Item(26, 5).Value = [double]"222.2702448282003"
$ws.Cells.Item(26, 6).Value = [double]"0.2771755285796327"
$ws.Cells.Item(26, 7).Value = [double]"6822.902463561311"
$ws.Cells.Item(26, 8).Value = [double]"0.4554659354423276"
$ws.Cells.Item(27, 2).Value = [double]"30.94477249000234"
$ws.Cells.Item(27, 3).Value = [double]"21.99999997884393"
$ws.Cells.Item(27, 4).Value = [double]"56"
$ws.Cells.Item(27, 5).Value = [double]"351.4163166375222"
$ws.Cells.Item(27, 6).Value = [double]"0.2893911295629741"
$ws.Cells.Item(27, 7).Value = [double]"6818.330678440493"
$ws.Cells.Item(27, 8).Value = [double]"0.4538467544240625"
$ws.Cells.Item(28, 2).Value = [double]"30.83062286000491"
$ws.Cells.Item(28, 3).Value = [double]"38"
$ws.Cells.Item(28, 4).Value = [double]"90"
$ws.Cells.Item(28, 5).Value = [double]"460.5643952549586"
$ws.Cells.Item(28, 6).Value = [double]"0.279819520828324"
$ws.Cells.Item(28, 7).Value = [double]"6878.909397507718"
$ws.Cells.Item(28, 8).Value = [double]"0.4481905644981323"
$ws.Cells.Item(29, 2).Value = [double]"30.72309198001023"
$ws.Cells.Item(29, 3).Value = [double]"33.99999971617449"
$ws.Cells.Item(29, 4).Value = [double]"110"
$ws.Cells.Item(29, 5).Value = [double]"445.7738139067373"
$ws.Cells.Item(29, 6).Value = [double]"0.2816342350166668"
$ws.Cells.Item(29, 7).Value = [double]"6987.707343726833"
$ws.Cells.Item(29, 8).Value = [double]"0.4396734217495769"
